# Insert a new worksheet "Destruction1" between "BaselineSoilSamples" and
# "SoilHealthSamples", populate it with depth-category labels, and make it
# the active sheet/tab.

$wb = $excel.ActiveWorkbook

# "SoilHealthSamples" is currently the 2nd sheet; add the new sheet before it
# so ordering becomes: BaselineSoilSamples, Destruction1, SoilHealthSamples.
$target = $wb.Worksheets.Item("SoilHealthSamples")
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "Destruction1"

$newSheet.Range("B1").Value = "Depth"
$newSheet.Range("B2").Value = "0-7.5"
$newSheet.Range("B3").Value = "7.5-15"
$newSheet.Range("B4").Value = "15-"

[void]$newSheet.Select()
[void]$newSheet.Range("A4").Select()
